$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8

# Columns A (date-looking, e.g. "2025-09-24") and C (all-digit, e.g.
# "250924") would be auto-converted by Excel's smart input parsing into a
# real date serial / number if assigned as a plain Value. Build each as a
# literal-string formula first, then convert the formula to a plain value
# in place (copy / paste-values) so the final cell holds a true text
# constant - matching the source file - without leaving a residual
# formula or a new number-format style behind.
$ws.Range("A$row").Formula = '="2025-09-24"'
$ws.Range("C$row").Formula = '="250924"'

$rngA = $ws.Range("A$row")
$rngA.Copy()
$rngA.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$rngC = $ws.Range("C$row")
$rngC.Copy()
$rngC.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

# Columns B, D, E are safe as plain text (Excel can't coerce "Pick 4",
# "3-2-8-1", or the timestamp string to a number/date), so a normal Value
# assignment keeps them literal text.
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 4).Value = "3-2-8-1"
$ws.Cells.Item($row, 5).Value = "2025-09-24T21:38:26.486+04:00"

$wb.Save()
